# introduce separate fee currency column, modify accordingly,
# make forex gains from dividends tax-free

$wb = $excel.ActiveWorkbook

# --- buy_orders: date,symbol,quantity,buy_price,fees,currency,comment
#     -> date,symbol,quantity,buy_price,currency,fees,fee_currency,comment
$ws = $wb.Worksheets.Item("buy_orders")
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "comment"
$ws.Range("G1").Value = "fee_currency"
$ws.Range("F1").Value = "fees"
$ws.Range("E1").Value = "currency"

# --- sell_orders: date,symbol,quantity,sell_price,fees,currency,comment
#     -> date,symbol,quantity,sell_price,currency,fees,fee_currency,comment
$ws = $wb.Worksheets.Item("sell_orders")
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "comment"
$ws.Range("G1").Value = "fee_currency"
$ws.Range("F1").Value = "fees"
$ws.Range("E1").Value = "currency"

# --- currency_conversions: date,source_amount,source_fees,source_currency,
#       target_amount,target_fees,target_currency,comment
#     -> date,source_amount,source_currency,target_amount,target_currency,
#       fees,fee_currency,comment
$ws = $wb.Worksheets.Item("currency_conversions")
$ws.Range("B1").Value = "source_amount"
$ws.Range("C1").Value = "source_currency"
$ws.Range("D1").Value = "target_amount"
$ws.Range("E1").Value = "target_currency"
$ws.Range("F1").Value = "fees"
$ws.Range("G1").Value = "fee_currency"
$ws.Range("H1").Value = "comment"

# --- money_transfers: date,buy_date,amount,fees,currency,comment
#     -> date,buy_date,amount,currency,fees,fee_currency,comment
$ws = $wb.Worksheets.Item("money_transfers")
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "comment"
$ws.Range("F1").Value = "fee_currency"
$ws.Range("E1").Value = "fees"
$ws.Range("D1").Value = "currency"

# the active tab moves from currency_conversions to sell_orders
$ws = $wb.Worksheets.Item("sell_orders")
$ws.Activate()
